$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B:E are treated as plain text so numeric-looking values
# (prices, percentages) are stored as strings, matching the source data.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "67.241.89"
$ws.Range("E2").Value = "  +6.20%  "

$ws.Range("D3").Value = "3.532.90"
$ws.Range("E3").Value = "  +8.47%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "193.06"
$ws.Range("E5").Value = "  +9.44%  "

$ws.Range("D6").Value = "558.85"
$ws.Range("E6").Value = "  +7.07%  "

$ws.Range("D7").Value = "3.525.83"
$ws.Range("E7").Value = "  +8.26%  "

$ws.Range("E8").Value = "  +2.89%  "

$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("D10").Value = "0.642"
$ws.Range("E10").Value = "  +6.35%  "

$ws.Range("E11").Value = "  +16.16%  "

$ws.Range("D12").Value = "56.63"
$ws.Range("E12").Value = "  +6.35%  "

$ws.Range("D13").Value = "0.0000273"
$ws.Range("E13").Value = "  +7.92%  "

$ws.Range("D14").Value = "9.53"
$ws.Range("E14").Value = "  +6.35%  "

$ws.Range("D15").Value = "4.107.41"
$ws.Range("E15").Value = "  +9.02%  "

$ws.Range("D16").Value = "3.533.67"
$ws.Range("E16").Value = "  +8.65%  "

$ws.Range("E17").Value = "  +5.34%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "18.44"
$ws.Range("E18").Value = "  +6.95%  "

$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "67.321.14"
$ws.Range("E19").Value = "  +6.42%  "

$ws.Range("D20").Value = "12.00"
$ws.Range("E20").Value = "  +8.55%  "

$ws.Range("E21").Value = "  +4.42%  "

$ws.Range("D22").Value = "409.61"
$ws.Range("E22").Value = "  +11.40%  "

$ws.Range("E23").Value = "  +6.51%  "

$ws.Range("D24").Value = "85.88"
$ws.Range("E24").Value = "  +6.60%  "

$ws.Range("D25").Value = "4.25"
$ws.Range("E25").Value = "  +8.79%  "

$ws.Range("D26").Value = "11.28"
$ws.Range("E26").Value = "  +1.92%  "

$ws.Range("E27").Value = "  +13.41%  "

$ws.Range("E28").Value = "  +0.67%  "

$ws.Range("D29").Value = "12.07"
$ws.Range("E29").Value = "  +7.08%  "

$ws.Range("D30").Value = "8.90"
$ws.Range("E30").Value = "  +8.62%  "

$ws.Range("D31").Value = "30.65"
$ws.Range("E31").Value = "  +7.89%  "

$ws.Range("D32").Value = "680.50"
$ws.Range("E32").Value = "  +3.80%  "

$ws.Range("D33").Value = "6.82"
$ws.Range("E33").Value = "  +6.75%  "

$ws.Range("D34").Value = "11.87"
$ws.Range("E34").Value = "  +6.38%  "

$ws.Range("E35").Value = "  +7.60%  "

$ws.Range("D36").Value = "60.58"

$ws.Range("D37").Value = "39.32"
$ws.Range("E37").Value = "  +7.34%  "

$ws.Range("D38").Value = "0.0₃0824"
$ws.Range("E38").Value = "  +14.86%  "

$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("D40").Value = "0.398"
$ws.Range("E40").Value = "  +5.72%  "

$ws.Range("E41").Value = "  +13.10%  "

$ws.Range("D42").Value = "3.40"
$ws.Range("E42").Value = "  +20.94%  "

$ws.Range("D43").Value = "3.03"
$ws.Range("E43").Value = "  +17.45%  "

$ws.Range("E44").Value = "  +0.33%  "

$ws.Range("D45").Value = "2.68"
$ws.Range("E45").Value = "  +6.87%  "

$ws.Range("D46").Value = "3.022.10"
$ws.Range("E46").Value = "  +4.34%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "3.32"
$ws.Range("E47").Value = "  +11.44%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "0.0421"
$ws.Range("E48").Value = "  +7.31%  "

$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "9.13"
$ws.Range("E49").Value = "  +18.12%  "

$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").Value = "2.74"
$ws.Range("E50").Value = "  +3.23%  "

$ws.Range("D51").Value = "0.131"
$ws.Range("E51").Value = "  +6.84%  "
